$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused explicit body style (xf index 4, "s=4") from the
# existing data rows (2-13) in columns A:C - these cells revert to the
# default (unstyled) cell format, matching the re-saved workbook.
$ws.Range("A2:C13").ClearFormats() | Out-Null

# New language blocks: kan, hin, tam - each repeating the same
# txt/xml/json/html x is_active=TRUE pattern used by eng/fra/ara.
$langs = @("kan", "hin", "tam")
$fileTypes = @(
    @("txt", "Text File"),
    @("xml", "XML File"),
    @("json", "Json File"),
    @("html", "html file")
)

$row = 14
foreach ($lang in $langs) {
    foreach ($ft in $fileTypes) {
        $ws.Range("A$row").Value = $lang
        $ws.Range("B$row").Value = $ft[0]
        $ws.Range("C$row").Value = $ft[1]

        # Copy D2 (style s="1", shared-string "TRUE") so the new is_active
        # cell is written as text "TRUE" (not an auto-coerced boolean) and
        # picks up the same number-format style as the rest of column D.
        $ws.Range("D2").Copy() | Out-Null
        $ws.Range("D$row").PasteSpecial() | Out-Null

        $row++
    }
}

$ws.Range("E24").Select() | Out-Null
